# Updated results and code
# Applies refreshed statistics to rows 4-13 (one row per metric) of the active sheet,
# matching the recomputed values from the updated analysis script.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: stats_for_precision
$ws.Range("B4").Value = 0.306
$ws.Range("D4").Value = 0.227
$ws.Range("E4").Value = 0.16
$ws.Range("G4").Value = 0.112
$ws.Range("H4").Value = 0.185
$ws.Range("J4").Value = 0.108
$ws.Range("K4").Value = 0.361
$ws.Range("L4").Value = 0.101
$ws.Range("M4").Value = 0.318
$ws.Range("N4").Value = 0.27
$ws.Range("O4").Value = 0.021
$ws.Range("P4").Value = 0.146
$ws.Range("Q4").Value = 0.528
$ws.Range("R4").Value = 0.221
$ws.Range("S4").Value = 0.471
$ws.Range("T4").Value = 0.271
$ws.Range("U4").Value = 0.09
$ws.Range("W4").Value = 0.239
$ws.Range("Y4").Value = 0.212
$ws.Range("Z4").Value = 0.439
$ws.Range("AA4").Value = 0.131
$ws.Range("AB4").Value = 0.362
$ws.Range("AE4").Value = 0.079
$ws.Range("AF4").Value = 0.757
$ws.Range("AG4").Value = 0.091
$ws.Range("AH4").Value = 0.301
$ws.Range("AI4").Value = 0.694
$ws.Range("AJ4").Value = 0.157
$ws.Range("AK4").Value = 0.396
$ws.Range("AL4").Value = 0.7
$ws.Range("AN4").Value = 0.343
$ws.Range("AO4").Value = 0.717

# Row 5: stats_for_recall
$ws.Range("B5").Value = 0.833
$ws.Range("C5").Value = 0.139
$ws.Range("D5").Value = 0.373
$ws.Range("E5").Value = 0.694
$ws.Range("F5").Value = 0.212
$ws.Range("G5").Value = 0.461
$ws.Range("H5").Value = 0.806
$ws.Range("I5").Value = 0.157
$ws.Range("J5").Value = 0.396
$ws.Range("K5").Value = 0.694
$ws.Range("L5").Value = 0.212
$ws.Range("M5").Value = 0.461
$ws.Range("N5").Value = 0.833
$ws.Range("O5").Value = 0.139
$ws.Range("P5").Value = 0.373
$ws.Range("Q5").Value = 0.583
$ws.Range("R5").Value = 0.243
$ws.Range("S5").Value = 0.493
$ws.Range("T5").Value = 0.556
$ws.Range("U5").Value = 0.247
$ws.Range("V5").Value = 0.497
$ws.Range("W5").Value = 0.722
$ws.Range("X5").Value = 0.201
$ws.Range("Y5").Value = 0.448
$ws.Range("Z5").Value = 0.806
$ws.Range("AA5").Value = 0.157
$ws.Range("AB5").Value = 0.396
$ws.Range("AC5").Value = 0.75
$ws.Range("AD5").Value = 0.188
$ws.Range("AE5").Value = 0.433
$ws.Range("AF5").Value = 0.972
$ws.Range("AG5").Value = 0.027
$ws.Range("AH5").Value = 0.164
$ws.Range("AI5").Value = 0.806
$ws.Range("AJ5").Value = 0.157
$ws.Range("AK5").Value = 0.396
$ws.Range("AL5").Value = 0.917
$ws.Range("AM5").Value = 0.076
$ws.Range("AN5").Value = 0.276
$ws.Range("AO5").Value = 0.898

# Row 6: stats_for_f1-score
$ws.Range("B6").Value = 0.448
$ws.Range("E6").Value = 0.26
$ws.Range("H6").Value = 0.301
$ws.Range("K6").Value = 0.475
$ws.Range("N6").Value = 0.408
$ws.Range("Q6").Value = 0.554
$ws.Range("T6").Value = 0.364
$ws.Range("W6").Value = 0.359
$ws.Range("Z6").Value = 0.5679999999999999
$ws.Range("AF6").Value = 0.851
$ws.Range("AI6").Value = 0.746
$ws.Range("AL6").Value = 0.794
$ws.Range("AO6").Value = 0.797

# Row 7: stats_for_f2-score
$ws.Range("B7").Value = 0.62
$ws.Range("E7").Value = 0.416
$ws.Range("H7").Value = 0.482
$ws.Range("K7").Value = 0.586
$ws.Range("N7").Value = 0.588
$ws.Range("Q7").Value = 0.571
$ws.Range("T7").Value = 0.459
$ws.Range("W7").Value = 0.514
$ws.Range("Z7").Value = 0.6909999999999999
$ws.Range("AC7").Value = 0.371
$ws.Range("AF7").Value = 0.92
$ws.Range("AI7").Value = 0.781
$ws.Range("AL7").Value = 0.863
$ws.Range("AO7").Value = 0.855

# Row 8: stats_for_NDCG
$ws.Range("B8").Value = 0.775
$ws.Range("C8").Value = 0.141
$ws.Range("D8").Value = 0.376
$ws.Range("E8").Value = 0.578
$ws.Range("H8").Value = 0.697
$ws.Range("I8").Value = 0.158
$ws.Range("J8").Value = 0.397
$ws.Range("K8").Value = 0.618
$ws.Range("L8").Value = 0.195
$ws.Range("M8").Value = 0.442
$ws.Range("N8").Value = 0.735
$ws.Range("O8").Value = 0.141
$ws.Range("P8").Value = 0.376
$ws.Range("Q8").Value = 0.5629999999999999
$ws.Range("R8").Value = 0.233
$ws.Range("S8").Value = 0.483
$ws.Range("T8").Value = 0.483
$ws.Range("V8").Value = 0.458
$ws.Range("W8").Value = 0.643
$ws.Range("X8").Value = 0.183
$ws.Range("Y8").Value = 0.428
$ws.Range("Z8").Value = 0.751
$ws.Range("AA8").Value = 0.154
$ws.Range("AB8").Value = 0.393
$ws.Range("AC8").Value = 0.636
$ws.Range("AD8").Value = 0.179
$ws.Range("AE8").Value = 0.424
$ws.Range("AF8").Value = 0.907
$ws.Range("AG8").Value = 0.045
$ws.Range("AH8").Value = 0.212
$ws.Range("AI8").Value = 0.795
$ws.Range("AJ8").Value = 0.156
$ws.Range("AK8").Value = 0.395
$ws.Range("AL8").Value = 0.886
$ws.Range("AM8").Value = 0.082
$ws.Range("AN8").Value = 0.286
$ws.Range("AO8").Value = 0.863

# Row 9: stats_for_M1
$ws.Range("B9").Value = 0.694
$ws.Range("C9").Value = 0.212
$ws.Range("D9").Value = 0.461
$ws.Range("E9").Value = 0.444
$ws.Range("F9").Value = 0.247
$ws.Range("G9").Value = 0.497
$ws.Range("H9").Value = 0.583
$ws.Range("I9").Value = 0.243
$ws.Range("J9").Value = 0.493
$ws.Range("K9").Value = 0.528
$ws.Range("L9").Value = 0.249
$ws.Range("M9").Value = 0.499
$ws.Range("N9").Value = 0.611
$ws.Range("O9").Value = 0.238
$ws.Range("P9").Value = 0.487
$ws.Range("Q9").Value = 0.528
$ws.Range("R9").Value = 0.249
$ws.Range("S9").Value = 0.499
$ws.Range("T9").Value = 0.389
$ws.Range("U9").Value = 0.238
$ws.Range("V9").Value = 0.487
$ws.Range("W9").Value = 0.528
$ws.Range("X9").Value = 0.249
$ws.Range("Y9").Value = 0.499
$ws.Range("Z9").Value = 0.667
$ws.Range("AA9").Value = 0.222
$ws.Range("AB9").Value = 0.471
$ws.Range("AC9").Value = 0.528
$ws.Range("AD9").Value = 0.249
$ws.Range("AE9").Value = 0.499
$ws.Range("AF9").Value = 0.806
$ws.Range("AG9").Value = 0.157
$ws.Range("AH9").Value = 0.396
$ws.Range("AI9").Value = 0.778
$ws.Range("AJ9").Value = 0.173
$ws.Range("AK9").Value = 0.416
$ws.Range("AL9").Value = 0.833
$ws.Range("AM9").Value = 0.139
$ws.Range("AN9").Value = 0.373
$ws.Range("AO9").Value = 0.806

# Row 10: stats_for_M3
$ws.Range("B10").Value = 0.833
$ws.Range("C10").Value = 0.139
$ws.Range("D10").Value = 0.373
$ws.Range("E10").Value = 0.611
$ws.Range("F10").Value = 0.238
$ws.Range("G10").Value = 0.487
$ws.Range("H10").Value = 0.722
$ws.Range("I10").Value = 0.201
$ws.Range("J10").Value = 0.448
$ws.Range("K10").Value = 0.694
$ws.Range("L10").Value = 0.212
$ws.Range("M10").Value = 0.461
$ws.Range("N10").Value = 0.806
$ws.Range("O10").Value = 0.157
$ws.Range("P10").Value = 0.396
$ws.Range("Q10").Value = 0.583
$ws.Range("R10").Value = 0.243
$ws.Range("S10").Value = 0.493
$ws.Range("T10").Value = 0.556
$ws.Range("U10").Value = 0.247
$ws.Range("V10").Value = 0.497
$ws.Range("W10").Value = 0.722
$ws.Range("X10").Value = 0.201
$ws.Range("Y10").Value = 0.448
$ws.Range("Z10").Value = 0.806
$ws.Range("AA10").Value = 0.157
$ws.Range("AB10").Value = 0.396
$ws.Range("AC10").Value = 0.639
$ws.Range("AD10").Value = 0.231
$ws.Range("AE10").Value = 0.48
$ws.Range("AF10").Value = 0.972
$ws.Range("AG10").Value = 0.027
$ws.Range("AH10").Value = 0.164
$ws.Range("AI10").Value = 0.806
$ws.Range("AJ10").Value = 0.157
$ws.Range("AK10").Value = 0.396
$ws.Range("AL10").Value = 0.917
$ws.Range("AM10").Value = 0.076
$ws.Range("AN10").Value = 0.276
$ws.Range("AO10").Value = 0.898

# Row 11: stats_for_M5
$ws.Range("B11").Value = 0.833
$ws.Range("C11").Value = 0.139
$ws.Range("D11").Value = 0.373
$ws.Range("E11").Value = 0.694
$ws.Range("F11").Value = 0.212
$ws.Range("G11").Value = 0.461
$ws.Range("H11").Value = 0.806
$ws.Range("I11").Value = 0.157
$ws.Range("J11").Value = 0.396
$ws.Range("K11").Value = 0.694
$ws.Range("L11").Value = 0.212
$ws.Range("M11").Value = 0.461
$ws.Range("N11").Value = 0.833
$ws.Range("O11").Value = 0.139
$ws.Range("P11").Value = 0.373
$ws.Range("Q11").Value = 0.583
$ws.Range("R11").Value = 0.243
$ws.Range("S11").Value = 0.493
$ws.Range("T11").Value = 0.556
$ws.Range("U11").Value = 0.247
$ws.Range("V11").Value = 0.497
$ws.Range("W11").Value = 0.722
$ws.Range("X11").Value = 0.201
$ws.Range("Y11").Value = 0.448
$ws.Range("Z11").Value = 0.806
$ws.Range("AA11").Value = 0.157
$ws.Range("AB11").Value = 0.396
$ws.Range("AC11").Value = 0.694
$ws.Range("AD11").Value = 0.212
$ws.Range("AE11").Value = 0.461
$ws.Range("AF11").Value = 0.972
$ws.Range("AG11").Value = 0.027
$ws.Range("AH11").Value = 0.164
$ws.Range("AI11").Value = 0.806
$ws.Range("AJ11").Value = 0.157
$ws.Range("AK11").Value = 0.396
$ws.Range("AL11").Value = 0.917
$ws.Range("AM11").Value = 0.076
$ws.Range("AN11").Value = 0.276
$ws.Range("AO11").Value = 0.898

# Row 12: stats_for_position
$ws.Range("B12").Value = 1.233
$ws.Range("C12").Value = 0.312
$ws.Range("D12").Value = 0.5590000000000001
$ws.Range("E12").Value = 1.68
$ws.Range("F12").Value = 1.098
$ws.Range("G12").Value = 1.048
$ws.Range("H12").Value = 1.621
$ws.Range("I12").Value = 1.408
$ws.Range("J12").Value = 1.187
$ws.Range("K12").Value = 1.4
$ws.Range("L12").Value = 0.5600000000000001
$ws.Range("M12").Value = 0.748
$ws.Range("N12").Value = 1.433
$ws.Range("O12").Value = 0.646
$ws.Range("P12").Value = 0.803
$ws.Range("Z12").Value = 1.207
$ws.Range("AA12").Value = 0.233
$ws.Range("AB12").Value = 0.483
$ws.Range("AC12").Value = 1.852
$ws.Range("AD12").Value = 2.571
$ws.Range("AE12").Value = 1.603
$ws.Range("AF12").Value = 1.2
$ws.Range("AG12").Value = 0.217
$ws.Range("AH12").Value = 0.466
$ws.Range("AI12").Value = 1.034
$ws.Range("AJ12").Value = 0.033
$ws.Range("AK12").Value = 0.182
$ws.Range("AL12").Value = 1.091
$ws.Range("AM12").Value = 0.083
$ws.Range("AN12").Value = 0.287
$ws.Range("AO12").Value = 1.108

# Row 13: stats_for_length (x of gs)
$ws.Range("B13").Value = 3.444
$ws.Range("C13").Value = 1.525
$ws.Range("D13").Value = 1.235
$ws.Range("E13").Value = 4.594
$ws.Range("F13").Value = 0.429
$ws.Range("G13").Value = 0.655
$ws.Range("H13").Value = 4.647
$ws.Range("I13").Value = 0.581
$ws.Range("J13").Value = 0.762
$ws.Range("K13").Value = 2.25
$ws.Range("L13").Value = 0.625
$ws.Range("M13").Value = 0.791
$ws.Range("N13").Value = 3.25
$ws.Range("O13").Value = 0.743
$ws.Range("P13").Value = 0.862
$ws.Range("Z13").Value = 2.545
$ws.Range("AA13").Value = 2.975
$ws.Range("AB13").Value = 1.725
$ws.Range("AC13").Value = 6.4
$ws.Range("AD13").Value = 2.24
$ws.Range("AE13").Value = 1.497
$ws.Range("AF13").Value = 1.528
$ws.Range("AG13").Value = 0.583
$ws.Range("AH13").Value = 0.763
$ws.Range("AI13").Value = 1.306
$ws.Range("AJ13").Value = 0.379
$ws.Range("AK13").Value = 0.616
$ws.Range("AL13").Value = 1.583
$ws.Range("AM13").Value = 0.743
$ws.Range("AN13").Value = 0.862
$ws.Range("AO13").Value = 1.472
